$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New columns AD/AE/AF hold the team's Wins/Losses/Ties record.
# Copy the header formatting (bold, bordered, centered) from the last
# existing header cell (AC1) onto the three new header cells.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Every player row (2-57) gets the same team record values.
for ($r = 2; $r -le 57; $r++) {
    $ws.Cells.Item($r, 30).Value = 80
    $ws.Cells.Item($r, 31).Value = 82
    $ws.Cells.Item($r, 32).Value = 0
}
